# Apply weekly update to the Berenjena (Macroferia Regional de Talca) sheet.
# The edit inserts three new daily price records into the historical table.
# Because the sheet is ordered with the newest-looking records shuffled in
# among older ones, the three new rows are inserted at the specific
# positions below; Excel's native row-insert shifts all subsequent rows
# down automatically, which reproduces the row renumbering seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the three new rows (top to bottom) --------------------------
# Inserting sequentially like this means each insertion point is expressed
# in terms of the row numbers that exist at the time of that particular
# insertion (i.e. after the previous inserts have already shifted things
# down).
$ws.Rows("8:8").Insert()
$ws.Rows("35:35").Insert()
$ws.Rows("50:50").Insert()

# --- Fill in the data for the new row 8 ----------------------------------
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44552
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112001
$ws.Range("G8").Value = "Berenjena"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("N8").Value = "$/caja 50 unidades"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 220
$ws.Range("Q8").Value = 50
$ws.Range("R8").Value = "Hortaliza"

# --- Fill in the data for the new row 35 ---------------------------------
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44554
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112001
$ws.Range("G35").Value = "Berenjena"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("N35").Value = "$/caja 50 unidades"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 200
$ws.Range("Q35").Value = 50
$ws.Range("R35").Value = "Hortaliza"

# --- Fill in the data for the new row 50 ---------------------------------
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44553
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = 100112001
$ws.Range("G50").Value = "Berenjena"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = 10000
$ws.Range("N50").Value = "$/caja 50 unidades"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 200
$ws.Range("Q50").Value = 50
$ws.Range("R50").Value = "Hortaliza"
